$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 23: current_phase D23 1 -> 2
$ws.Range("D23").Value = 2

# Row 24: current_phase D24 1 -> 2
$ws.Range("D24").Value = 2

# Row 24: last_action_date E24 updated timestamp (keep as text)
$ws.Range("E24").Value = "2026-02-12T15:39:45.630203+00:00"

# Row 24: reactions_count H24 3 -> 5
$ws.Range("H24").Value = 5

# Row 24: reacted_message_ids L24 append new ids
$ws.Range("L24").Value = "[63, 31910, 19424, 30964, 30729]"
